# Issue #44360 - EPPlus can not calculate with Ranges in an IF-Function.
# Adds a small "which due-date quarter is this" lookup block (H1:K6) driven
# by four new single-cell named ranges (sn_duedate..sn_duedate3) and three
# helper ranges (IB_Q1_AJ..IB_Q3_AJ), plus a couple of incidental smoke-test
# cells/rows exercising the shared formulas already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New workbook-level defined names -------------------------------------
$wb.Names.Add("IB_Q1_AJ", "=Tabelle1!`$C`$2:`$C`$5")
$wb.Names.Add("IB_Q2_AJ", "=Tabelle1!`$B`$2:`$B`$5")
$wb.Names.Add("IB_Q3_AJ", "=Tabelle1!`$B`$2:`$B`$5")
$wb.Names.Add("sn_duedate", "=Tabelle1!`$E`$2")
$wb.Names.Add("sn_duedate1", "=Tabelle1!`$E`$3")
$wb.Names.Add("sn_duedate2", "=Tabelle1!`$E`$4")
$wb.Names.Add("sn_duedate3", "=Tabelle1!`$E`$5")

# --- E2:E5: text-formatted quarter-end labels referenced by sn_duedate* ---
$ws.Range("E2:E5").NumberFormat = "@"
$ws.Range("E2").Value = "31.3."
$ws.Range("E3").Value = "30.6."
$ws.Range("E4").Value = "30.9."
$ws.Range("E5").Value = "20.1."

# --- H1:K6: one nested-IF lookup per sn_duedate[n] name --------------------
$dueNames = @("sn_duedate", "sn_duedate1", "sn_duedate2", "sn_duedate3")
$cols = @("H", "I", "J", "K")
for ($c = 0; $c -lt 4; $c++) {
    $col = $cols[$c]
    $nm = $dueNames[$c]
    $formula = '=IF(' + $nm + '="31.3.",IB_Q1_AJ,IF(' + $nm + '="30.6.",IB_Q2_AJ,IF(' + $nm + '="30.9.",IB_Q3_AJ,"Falsche Auswahl")))'
    for ($r = 1; $r -le 6; $r++) {
        $ws.Range($col + $r).Formula = $formula
    }
}

# Column K is a bit wider so the lookup results aren't truncated.
$ws.Columns.Item(11).ColumnWidth = 10.26

# --- F11 / F12: small IF() smoke-test formulas ------------------------------
$ws.Range("F11").Formula = "=IF(0=1,1,2)"
$ws.Range("F12").Formula = "=IF(0=0,1,2)"

# --- Row 16: extra data row, extends the range1+range2 shared formula -----
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 20
$ws.Range("C16").Formula = "=range1 + range2"

# --- Row 21: extend the range3+range4 shared formula by one column on
#     each side (E21 before F21, S21 after R21) -----------------------------
$ws.Range("E21").Formula = "=range3+range4"
$ws.Range("S21").Formula = "=range3+range4"

# --- Sheet print setup / selection, matching the final authored state -----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

[void]$ws.Range("F14").Select()
